$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Finalized EchoPlay chip replacement planning: update the quantity counts.
$ws.Range("A5").Value = 2
$ws.Range("A6").Value = 1

# Leave the selection where the user ended up after finishing the edit.
$ws.Range("J8").Select()
